$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.384480357170105
$ws.Range("B1").Value = 2.638251304626465
$ws.Range("C1").Value = 5.990755081176758
$ws.Range("D1").Value = 2.316596031188965
$ws.Range("E1").Value = 1.20943820476532
